# Delete row 268 (XNN / Xining, China) from the data table.
# This shifts all subsequent rows up by one, which matches the diff:
# the old row 269 (FRU/Bishkek) becomes the new row 268, and so on,
# until the old last row 333 (YWG/Winnipeg) disappears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(268).Delete()
